# Add Wins/Losses/Ties columns (AD, AE, AF) to the season-record worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy the formatting from the last existing header cell (AC1)
# so the new header cells (AD1:AF1) get the same bold/bordered/centered style.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows: every player row (2-51) gets the team's season record.
$wins = 69
$losses = 93
$ties = 0

for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
